$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 5 data
$ws.Range("A5").Value = 4
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "1.3"
$ws.Range("C5").Value = "Dit uur was iets minder prodcutief. Ik heb de UI gebonden aan het generate script, en er gelijk voor gezorgd dat de user pas een doolhoof kan genereren als alle input velden een getal hebben. de functie om een doolhof te maken door op ""G"" te klikken is verwijderd."

# Copy styles from row 4 to row 5 for columns B and C
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)

# Update selection
$ws.Range("C19").Select()
